$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.934894
$ws.Range("H2").Value = 32.804682
$ws.Range("I2").Value = 0.3698068269583527
$ws.Range("J2").Value = 0.3698068269583527
$ws.Range("Q2").Value = 36.77125647906533
$ws.Range("R2").Value = 330.941308311588
$ws.Range("S2").Value = 0.2769442156917211
$ws.Range("T2").Value = 0.2769442156917211

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.934894
$ws.Range("H3").Value = 32.804682
$ws.Range("I3").Value = 0.3698068269583527
$ws.Range("J3").Value = 0.3698068269583527
$ws.Range("O3").Value = 0.251111132886386
$ws.Range("Q3").Value = 12.32982926786267
$ws.Range("R3").Value = 110.968463410764
$ws.Range("S3").Value = 0.09286261126663165
$ws.Range("T3").Value = 0.09286261126663163

# Row 4
$ws.Range("I4").Value = 0.3872921463699351
$ws.Range("J4").Value = 0.3872921463699351
$ws.Range("Q4").Value = 38.509886265838
$ws.Range("S4").Value = 0.2900387767369808
$ws.Range("T4").Value = 0.2900387767369807

# Row 5
$ws.Range("I5").Value = 0.3872921463699351
$ws.Range("J5").Value = 0.3872921463699351
$ws.Range("O5").Value = 0.251111132886386
$ws.Range("S5").Value = 0.09725336963295443
$ws.Range("T5").Value = 0.09725336963295442

# Row 6
$ws.Range("I6").Value = 0.2429010266717122
$ws.Range("J6").Value = 0.2429010266717122
$ws.Range("S6").Value = 0.1819058746849123
$ws.Range("T6").Value = 0.1819058746849123

# Row 7
$ws.Range("I7").Value = 0.2429010266717122
$ws.Range("J7").Value = 0.2429010266717122
$ws.Range("O7").Value = 0.251111132886386
$ws.Range("Q7").Value = 8.098628715115778
$ws.Range("R7").Value = 72.887658436042
$ws.Range("S7").Value = 0.0609951519867999
$ws.Range("T7").Value = 0.06099515198679988
